# Nov 29th - Status
# Append the daily status rows for 27/11, 28/11 (holidays) and 29/11/2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 136: 27/11/2021 - Holiday
$ws.Range("A136").Value = "27/11/2021"
$ws.Range("B136").Value = "HOLIDAY"

# Row 137: 28/11/2021 - Holiday
$ws.Range("A137").Value = "28/11/2021"
$ws.Range("B137").Value = "HOLIDAY"

# Row 138: 29/11/2021
$ws.Range("A138").Value = "29/11/2021"
$ws.Range("B138").Value = "Gone through android codes shared"
$ws.Range("C138").Value = "Revising basic CPP and OOPs concepts"

# Row 139: continuation of 29/11/2021 notes
$ws.Range("B139").Value = "Tried Logical and Puzzle questions"
$ws.Range("C139").Value = "Preparing for the interview"

[void]$ws.Range("B139").Select()
